# Update "想去人数" (F column) figures across the sheets.
# 全部类型 (All Types) is an aggregate sheet that mirrors rows from
# 展览 / 演出 / 本地生活, so the same events get updated there too.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsAll = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$wsExhibition.Range("F4").Value = 228
$wsExhibition.Range("F6").Value = 1139
$wsExhibition.Range("F7").Value = 917
$wsExhibition.Range("F11").Value = 885
$wsExhibition.Range("F12").Value = 311
$wsExhibition.Range("F17").Value = 1256
$wsExhibition.Range("F19").Value = 67
$wsExhibition.Range("F20").Value = 1537
$wsExhibition.Range("F21").Value = 1296
$wsExhibition.Range("F22").Value = 747
$wsExhibition.Range("F28").Value = 3270
$wsExhibition.Range("F29").Value = 639
$wsExhibition.Range("F31").Value = 1456

# 演出 (sheet2)
$wsPerformance.Range("F8").Value = 9
$wsPerformance.Range("F9").Value = 36

# 本地生活 (sheet3)
$wsLocalLife.Range("F2").Value = 769

# 全部类型 (sheet4) - aggregate of all sheets above
$wsAll.Range("F3").Value = 769
$wsAll.Range("F7").Value = 228
$wsAll.Range("F10").Value = 1139
$wsAll.Range("F11").Value = 917
$wsAll.Range("F18").Value = 9
$wsAll.Range("F19").Value = 36
$wsAll.Range("F22").Value = 885
$wsAll.Range("F23").Value = 311
$wsAll.Range("F28").Value = 1256
$wsAll.Range("F30").Value = 67
$wsAll.Range("F31").Value = 1537
$wsAll.Range("F32").Value = 1296
$wsAll.Range("F33").Value = 747
$wsAll.Range("F41").Value = 3270
$wsAll.Range("F42").Value = 639
$wsAll.Range("F44").Value = 1456
